$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text so numeric-looking values are not
# reinterpreted as numbers by Excel (e.g. "215.23" -> 215.23 float).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '26.825.20'
$ws.Range("E2").Value = '  -1.28%  '

# Row 3
$ws.Range("D3").Value = '1.662.14'
$ws.Range("E3").Value = '  +0.03%  '

# Row 4
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$ws.Range("D5").Value = '215.23'
$ws.Range("E5").Value = '  -0.22%  '

# Row 6
$ws.Range("D6").Value = '0.534'
$ws.Range("E6").Value = '  +5.06%  '

# Row 7
$ws.Range("E7").Value = '  +0.18%  '

# Row 8
$ws.Range("D8").Value = '0.250'
$ws.Range("E8").Value = '  +0.22%  '

# Row 9
$ws.Range("E9").Value = '  +0.38%  '

# Row 10
$ws.Range("D10").Value = '20.14'
$ws.Range("E10").Value = '  +2.20%  '

# Row 11
$ws.Range("D11").Value = '0.0896'
$ws.Range("E11").Value = '  +3.72%  '

# Row 12
$ws.Range("D12").Value = '1.898.72'
$ws.Range("E12").Value = '  +0.21%  '

# Row 13
$ws.Range("D13").Value = '1.676.63'
$ws.Range("E13").Value = '  +0.56%  '

# Row 14
$ws.Range("E14").Value = '  -0.20%  '

# Row 15
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '65.96'
$ws.Range("E15").Value = '  +1.49%  '

# Row 16
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").Value = '0.522'
$ws.Range("E16").Value = '  +0.06%  '

# Row 17
$ws.Range("D17").Value = '26.856.83'
$ws.Range("E17").Value = '  -1.07%  '

# Row 18
$ws.Range("D18").Value = '231.27'
$ws.Range("E18").Value = '  -4.04%  '

# Row 19
$ws.Range("D19").Value = '7.79'
$ws.Range("E19").Value = '  -1.10%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0729'
$ws.Range("E20").Value = '  -0.22%  '

# Row 21
$ws.Range("E21").Value = '  +0.08%  '

# Row 22
$ws.Range("D22").Value = '4.43'
$ws.Range("E22").Value = '  -0.47%  '

# Row 23
$ws.Range("E23").Value = '  -2.26%  '

# Row 24
$ws.Range("D24").Value = '9.14'
$ws.Range("E24").Value = '  -2.14%  '

# Row 25
$ws.Range("D25").Value = '145.64'
$ws.Range("E25").Value = '  -0.27%  '

# Row 26
$ws.Range("D26").Value = '7.11'
$ws.Range("E26").Value = '  -0.83%  '

# Row 27
$ws.Range("E27").Value = '  +1.21%  '

# Row 28
$ws.Range("D28").Value = '15.84'
$ws.Range("E28").Value = '  -0.28%  '

# Row 29
$ws.Range("E29").Value = '  +0.01%  '

# Row 30
$ws.Range("E30").Value = '  -0.68%  '

# Row 31
$ws.Range("E31").Value = '  -0.04%  '

# Row 32
$ws.Range("D32").Value = '3.35'
$ws.Range("E32").Value = '  +1.69%  '

# Row 33
$ws.Range("D33").Value = '1.460.50'
$ws.Range("E33").Value = '  -5.13%  '

# Row 34
$ws.Range("D34").Value = '3.14'
$ws.Range("E34").Value = '  +2.85%  '

# Row 35
$ws.Range("E35").Value = '  +2.33%  '

# Row 36
$ws.Range("D36").Value = '2.42'
$ws.Range("E36").Value = '  -0.23%  '

# Row 37
$ws.Range("D37").Value = '0.572'
$ws.Range("E37").Value = '  -0.86%  '

# Row 38
$ws.Range("D38").Value = '0.895'
$ws.Range("E38").Value = '  -0.09%  '

# Row 39
$ws.Range("D39").Value = '0.0168'
$ws.Range("E39").Value = '  -0.49%  '

# Row 40
$ws.Range("D40").Value = '5.87'
$ws.Range("E40").Value = '  -1.70%  '

# Row 41
$ws.Range("E41").Value = '  +0.13%  '

# Row 42
$ws.Range("D42").Value = '2.26'
$ws.Range("E42").Value = '  -0.96%  '

# Row 43
$ws.Range("D43").Value = '0.978'
$ws.Range("E43").Value = '  +6.07%  '

# Row 44
$ws.Range("D44").Value = '65.65'
$ws.Range("E44").Value = '  -0.95%  '

# Row 45
$ws.Range("D45").Value = '1.810.09'
$ws.Range("E45").Value = '  +0.50%  '

# Row 46
$ws.Range("D46").Value = '0.776'
$ws.Range("E46").Value = '  +0.37%  '

# Row 47
$ws.Range("D47").Value = '90.23'
$ws.Range("E47").Value = '  -0.19%  '

# Row 48
$ws.Range("E48").Value = '  -0.92%  '

# Row 49
$ws.Range("E49").Value = '  -1.60%  '

# Row 50
$ws.Range("E50").Value = '  +2.89%  '

# Row 51
$ws.Range("D51").Value = '0.0506'
$ws.Range("E51").Value = '  +0.34%  '
